# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to match the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 517
$ws1.Range("F3").Value  = 734
$ws1.Range("F4").Value  = 1443
$ws1.Range("F5").Value  = 220
$ws1.Range("F8").Value  = 6127
$ws1.Range("F10").Value = 397
$ws1.Range("F11").Value = 110
$ws1.Range("F12").Value = 4988
$ws1.Range("F15").Value = 1167
$ws1.Range("F18").Value = 57
$ws1.Range("F21").Value = 23
$ws1.Range("F22").Value = 3483
$ws1.Range("F23").Value = 145

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 517
$ws4.Range("F4").Value  = 734
$ws4.Range("F5").Value  = 1443
$ws4.Range("F6").Value  = 220
$ws4.Range("F9").Value  = 6127
$ws4.Range("F11").Value = 397
$ws4.Range("F12").Value = 110
$ws4.Range("F13").Value = 4988
$ws4.Range("F16").Value = 1167
$ws4.Range("F19").Value = 57
$ws4.Range("F22").Value = 23
$ws4.Range("F23").Value = 3483
$ws4.Range("F25").Value = 145
